# Bug fix, output image and output table headers.
# Rename the "Type" column values in A4:A7 from the old "M-..." labels to
# the new "[M-H]-..." labels, and set the worksheet page setup (paper size
# A4 / portrait) so printing produces the expected output image.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "[M-H]-sn2"
$ws.Range("A4").Value = "[M-H]-sn1"
$ws.Range("A6").Value = "[M-H]-sn1-H2O"
$ws.Range("A7").Value = "[M-H]-sn2-H2O"

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
